# bug fix in scale gblems_expanded
#
# Row 3 ("law_order") was missing the min/max scale-anchor labels that the
# other rows already have (perceived_intent uses "Extremely unlikely" /
# "Extremely likely"; serve_protect uses "Not at all" / "Extremely"). This
# fills in the missing "Not at all" / "Extremely" labels for law_order so it
# matches the serve_protect scale, and leaves the selection on the
# newly-fixed cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Not at all"
$ws.Range("D3").Value = "Extremely"

# Reflect the fix in the current selection/active cell.
$excel.Goto($ws.Range("C3:D4"))
